$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.280.16"
$ws.Range("E2").Value = "  -3.14%  "
$ws.Range("D3").Value = "3.171.09"
$ws.Range("E3").Value = "  -8.35%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.63"
$ws.Range("E5").Value = "  -4.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.77"
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "3.167.03"
$ws.Range("E9").Value = "  -8.45%  "
$ws.Range("E10").Value = "  -6.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.63"
$ws.Range("E11").Value = "  -4.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.397"
$ws.Range("E12").Value = "  -4.80%  "
$ws.Range("D13").Value = "3.712.56"
$ws.Range("E13").Value = "  -8.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.135"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.40"
$ws.Range("E15").Value = "  -7.47%  "
$ws.Range("D16").Value = "64.226.42"
$ws.Range("E16").Value = "  -3.07%  "
$ws.Range("E17").Value = "  -5.54%  "
$ws.Range("D18").Value = "3.164.41"
$ws.Range("E18").Value = "  -8.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.74"
$ws.Range("E20").Value = "  -5.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "354.30"
$ws.Range("E21").Value = "  -5.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.21"
$ws.Range("E22").Value = "  -5.49%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.21"
$ws.Range("E24").Value = "  -5.61%  "
$ws.Range("E25").Value = "  -6.28%  "
$ws.Range("E26").Value = "  -5.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.59"
$ws.Range("E27").Value = "  -2.80%  "
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.63"
$ws.Range("E30").Value = "  -4.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E33").Value = "  -7.02%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.21"
$ws.Range("E34").Value = "  -5.64%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.64"
$ws.Range("E35").Value = "  -5.73%  "
$ws.Range("E36").Value = "  -7.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "154.97"
$ws.Range("E37").Value = "  -4.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.811"
$ws.Range("E38").Value = "  -8.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.85"
$ws.Range("E39").Value = "  -9.51%  "
$ws.Range("E40").Value = "  -3.65%  "
$ws.Range("E41").Value = "  -6.23%  "
$ws.Range("D42").Value = "2.606.77"
$ws.Range("E42").Value = "  -6.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.19"
$ws.Range("E43").Value = "  -7.19%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.71"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.01"
$ws.Range("E45").Value = "  -6.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0660"
$ws.Range("E46").Value = "  -4.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.97"
$ws.Range("E47").Value = "  -5.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "321.90"
$ws.Range("E48").Value = "  -4.75%  "
$ws.Range("E49").Value = "  -7.48%  "
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("E51").Value = "  -0.11%  "
